# Update participant/attendee counts (column F) on several rows across sheets.
# Values below were derived from the canonical OOXML diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: "展览"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(5, 6).Value  = 271
$ws1.Cells.Item(7, 6).Value  = 836
$ws1.Cells.Item(8, 6).Value  = 26
$ws1.Cells.Item(9, 6).Value  = 485
$ws1.Cells.Item(10, 6).Value = 69
$ws1.Cells.Item(12, 6).Value = 102
$ws1.Cells.Item(14, 6).Value = 21
$ws1.Cells.Item(15, 6).Value = 402
$ws1.Cells.Item(16, 6).Value = 6519
$ws1.Cells.Item(20, 6).Value = 7448
$ws1.Cells.Item(22, 6).Value = 33
$ws1.Cells.Item(23, 6).Value = 3356
$ws1.Cells.Item(25, 6).Value = 1124
$ws1.Cells.Item(28, 6).Value = 345
$ws1.Cells.Item(29, 6).Value = 59
$ws1.Cells.Item(31, 6).Value = 180
$ws1.Cells.Item(32, 6).Value = 1518
$ws1.Cells.Item(37, 6).Value = 1135
$ws1.Cells.Item(38, 6).Value = 1634

# Sheet 2: "演出"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 65

# Sheet 3: "本地生活"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 6).Value = 1203

# Sheet 4: "全部类型"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(4, 6).Value  = 1203
$ws4.Cells.Item(8, 6).Value  = 271
$ws4.Cells.Item(10, 6).Value = 836
$ws4.Cells.Item(11, 6).Value = 26
$ws4.Cells.Item(12, 6).Value = 485
$ws4.Cells.Item(14, 6).Value = 69
$ws4.Cells.Item(16, 6).Value = 65
$ws4.Cells.Item(17, 6).Value = 102
$ws4.Cells.Item(19, 6).Value = 21
$ws4.Cells.Item(20, 6).Value = 402
$ws4.Cells.Item(21, 6).Value = 6519
$ws4.Cells.Item(25, 6).Value = 7448
$ws4.Cells.Item(27, 6).Value = 33
$ws4.Cells.Item(28, 6).Value = 3356
$ws4.Cells.Item(30, 6).Value = 1124
$ws4.Cells.Item(33, 6).Value = 345
$ws4.Cells.Item(34, 6).Value = 59
$ws4.Cells.Item(37, 6).Value = 180
$ws4.Cells.Item(38, 6).Value = 1518
$ws4.Cells.Item(43, 6).Value = 1135
$ws4.Cells.Item(44, 6).Value = 1634
